# "date time fixed using excel for han data consolidation file"
#
# Rows 23:40 of column A had been accidentally filled with a run of
# consecutive daily dates (2019-01-01 .. 2019-01-18, serials 43466-43483)
# using a "d-mmm" number format. They should instead continue the yearly
# "Jan 1" sequence established by rows 2:22 (1980..2000) through
# 2001..2018, formatted the same way as those earlier rows ("mmm-yy").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wp")

# Corrected Jan-1 serial date values for 2001 through 2018.
$correctedDates = @{
    23 = 36892   # 2001-01-01
    24 = 37257   # 2002-01-01
    25 = 37622   # 2003-01-01
    26 = 37987   # 2004-01-01
    27 = 38353   # 2005-01-01
    28 = 38718   # 2006-01-01
    29 = 39083   # 2007-01-01
    30 = 39448   # 2008-01-01
    31 = 39814   # 2009-01-01
    32 = 40179   # 2010-01-01
    33 = 40544   # 2011-01-01
    34 = 40909   # 2012-01-01
    35 = 41275   # 2013-01-01
    36 = 41640   # 2014-01-01
    37 = 42005   # 2015-01-01
    38 = 42370   # 2016-01-01
    39 = 42736   # 2017-01-01
    40 = 43101   # 2018-01-01
}

foreach ($row in $correctedDates.Keys) {
    $ws.Cells.Item($row, 1).Value = $correctedDates[$row]
}

# Match the number format used by the rest of the date column (rows 2:22),
# which collapses these cells back onto the same style record the others
# use instead of the stray "d-mmm" style.
$ws.Range("A23:A40").NumberFormat = "mmm-yy"

# The sheet was left scrolled down with A23:A40 selected (active cell A23).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 20
$win.ScrollColumn = 1
$ws.Range("A23:A40").Select()
